$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 100.8373843333333
$ws.Range("H2").Value = 302.512153
$ws.Range("I2").Value = 0.6551985585448407
$ws.Range("J2").Value = 0.6551985585448408
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.228108666666667
$ws.Range("N2").Value = 6.684326
$ws.Range("O2").Value = 0.5687849952918405
$ws.Range("P2").Value = 0.5687849952918405
$ws.Range("Q2").Value = 224.6766499570976
$ws.Range("R2").Value = 2022.089849613878
$ws.Range("S2").Value = 0.3726671090371478
$ws.Range("T2").Value = 0.3726671090371479

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 100.8373843333333
$ws.Range("H3").Value = 302.512153
$ws.Range("I3").Value = 0.6551985585448407
$ws.Range("J3").Value = 0.6551985585448408
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.1124773333333333
$ws.Range("N3").Value = 0.337432
$ws.Range("O3").Value = 0.02871288122861097
$ws.Range("P3").Value = 0.02871288122861097
$ws.Range("Q3").Value = 11.34192009012178
$ws.Range("R3").Value = 102.077280811096
$ws.Range("S3").Value = 0.01881263839265512
$ws.Range("T3").Value = 0.01881263839265513

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 100.8373843333333
$ws.Range("H4").Value = 302.512153
$ws.Range("I4").Value = 0.6551985585448407
$ws.Range("J4").Value = 0.6551985585448408
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.576726666666667
$ws.Range("N4").Value = 4.73018
$ws.Range("O4").Value = 0.4025021234795486
$ws.Range("P4").Value = 0.4025021234795487
$ws.Range("Q4").Value = 158.9929928752822
$ws.Range("R4").Value = 1430.93693587754
$ws.Range("S4").Value = 0.2637188111150377
$ws.Range("T4").Value = 0.2637188111150378

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 23.90796933333333
$ws.Range("H5").Value = 71.723908
$ws.Range("I5").Value = 0.1553438454249564
$ws.Range("J5").Value = 0.1553438454249564
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.228108666666667
$ws.Range("N5").Value = 6.684326
$ws.Range("O5").Value = 0.5687849952918405
$ws.Range("P5").Value = 0.5687849952918405
$ws.Range("Q5").Value = 53.26955367400089
$ws.Range("R5").Value = 479.425983066008
$ws.Range("S5").Value = 0.08835724838865022
$ws.Range("T5").Value = 0.08835724838865022

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 23.90796933333333
$ws.Range("H6").Value = 71.723908
$ws.Range("I6").Value = 0.1553438454249564
$ws.Range("J6").Value = 0.1553438454249564
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.1124773333333333
$ws.Range("N6").Value = 0.337432
$ws.Range("O6").Value = 0.02871288122861097
$ws.Range("P6").Value = 0.02871288122861097
$ws.Range("Q6").Value = 2.689104636028444
$ws.Range("R6").Value = 24.201941724256
$ws.Range("S6").Value = 0.004460369383282475
$ws.Range("T6").Value = 0.004460369383282476

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 23.90796933333333
$ws.Range("H7").Value = 71.723908
$ws.Range("I7").Value = 0.1553438454249564
$ws.Range("J7").Value = 0.1553438454249564
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.576726666666667
$ws.Range("N7").Value = 4.73018
$ws.Range("O7").Value = 0.4025021234795486
$ws.Range("P7").Value = 0.4025021234795487
$ws.Range("Q7").Value = 37.69633279371555
$ws.Range("R7").Value = 339.2669951434399
$ws.Range("S7").Value = 0.06252622765302372
$ws.Range("T7").Value = 0.06252622765302372

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 29.15819666666667
$ws.Range("H8").Value = 87.47459
$ws.Range("I8").Value = 0.1894575960302029
$ws.Range("J8").Value = 0.1894575960302029
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.228108666666667
$ws.Range("N8").Value = 6.684326
$ws.Range("O8").Value = 0.5687849952918405
$ws.Range("P8").Value = 0.5687849952918405
$ws.Range("Q8").Value = 64.96763069737112
$ws.Range("R8").Value = 584.70867627634
$ws.Range("S8").Value = 0.1077606378660424
$ws.Range("T8").Value = 0.1077606378660424

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 29.15819666666667
$ws.Range("H9").Value = 87.47459
$ws.Range("I9").Value = 0.1894575960302029
$ws.Range("J9").Value = 0.1894575960302029
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.1124773333333333
$ws.Range("N9").Value = 0.337432
$ws.Range("O9").Value = 0.02871288122861097
$ws.Range("P9").Value = 0.02871288122861097
$ws.Range("Q9").Value = 3.279636205875555
$ws.Range("R9").Value = 29.51672585288
$ws.Range("S9").Value = 0.005439873452673374
$ws.Range("T9").Value = 0.005439873452673374

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 29.15819666666667
$ws.Range("H10").Value = 87.47459
$ws.Range("I10").Value = 0.1894575960302029
$ws.Range("J10").Value = 0.1894575960302029
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.576726666666667
$ws.Range("N10").Value = 4.73018
$ws.Range("O10").Value = 0.4025021234795486
$ws.Range("P10").Value = 0.4025021234795487
$ws.Range("Q10").Value = 45.97450623624444
$ws.Range("R10").Value = 413.7705561262
$ws.Range("S10").Value = 0.07625708471148718
$ws.Range("T10").Value = 0.07625708471148719
